$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the client location (B3) to the new project path
$ws.Range("B3").Value = "W:\Projects\בהת\135_סטריפ_דרך_חברון\קבצי עבודה\תחזיות_דמוגרפיות"

# Update the version date (B5) to the new value
$ws.Range("B5").Value = 240530

# Update the active selection to B6
$ws.Range("B6").Select()
